$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 37 (pushing existing rows down), matching the
# style/format of the surrounding data rows automatically.
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()

# Populate the two newly inserted rows with the new API catalogue entries.
$ws.Cells.Item(37, 2).Value = "transaction.delete.master.setBusinessDocumentNumbering"
$ws.Cells.Item(37, 3).Value = "Menghapus Data Penomoran Dokumen Bisnis"
$ws.Cells.Item(38, 2).Value = "transaction.delete.master.setBusinessDocumentNumberingFormat"
$ws.Cells.Item(38, 3).Value = "Menghapus Data Format Penomoran Dokumen Bisnis"
